$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.873.80'
$ws.Range('E2').Value = '  +4.10%  '
$ws.Range('D3').Value = '2.280.13'
$ws.Range('E3').Value = '  +4.77%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = "'250.78"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('E6').Value = '  +3.25%  '
$ws.Range('D7').Value = "'71.70"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +8.95%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = "'0.661"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +16.06%  '
$ws.Range('D10').Value = "'38.93"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.59%  '
$ws.Range('D11').Value = "'59.94"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.74%  '
$ws.Range('D12').Value = "'0.0973"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.90%  '
$ws.Range('D13').Value = "'7.42"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.16%  '
$ws.Range('D14').Value = "'0.105"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').Value = '2.621.54'
$ws.Range('E15').Value = '  +4.95%  '
$ws.Range('D16').Value = "'14.88"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.71%  '
$ws.Range('D17').Value = "'0.889"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.40%  '
$ws.Range('D18').Value = '2.275.02'
$ws.Range('E18').Value = '  +4.14%  '
$ws.Range('D19').Value = '42.815.78'
$ws.Range('E19').Value = '  +4.07%  '
$ws.Range('D20').Value = "'0.0000102"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.20%  '
$ws.Range('E21').Value = '  +3.97%  '
$ws.Range('D22').Value = "'73.25"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.26%  '
$ws.Range('D23').Value = "'236.20"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.69%  '
$ws.Range('E24').Value = '  +5.30%  '
$ws.Range('D25').Value = "'4.04"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.54%  '
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').Value = "'11.41"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.15%  '
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('D30').Value = "'2.13"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('D31').Value = "'167.87"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('D32').Value = "'21.04"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.17%  '
$ws.Range('D33').Value = "'6.47"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +14.28%  '
$ws.Range('E34').Value = '  +4.59%  '
$ws.Range('E35').Value = '  +9.02%  '
$ws.Range('D36').Value = "'31.42"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +29.03%  '
$ws.Range('E37').Value = '  +4.05%  '
$ws.Range('D38').Value = "'4.52"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +13.69%  '
$ws.Range('E39').Value = '  +5.30%  '
$ws.Range('D40').Value = "'0.0314"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.98%  '
$ws.Range('D41').Value = "'2.33"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.05%  '
$ws.Range('D42').Value = "'13.30"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +18.34%  '
$ws.Range('D43').Value = "'5.84"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.81%  '
$ws.Range('D44').Value = "'0.209"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.97%  '
$ws.Range('D45').Value = "'5.02"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.09%  '
$ws.Range('E46').Value = '  +8.00%  '
$ws.Range('D47').Value = "'62.24"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.70%  '
$ws.Range('E48').Value = '  +2.63%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  +3.24%  '
$ws.Range('E51').Value = '  +4.76%  '
